# Refresh crypto snapshot values (price / 1h volume change) scraped from
# coinranking.com, matching the GitHub Actions automated update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.025.64'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.79%  '
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.423.53'
$ws.Range("D3").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.91'
$ws.Range("D5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.57'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.01%  '
$ws.Range("E6").NumberFormat = "General"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E7").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("E8").NumberFormat = "General"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.04'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.83%  '
$ws.Range("E9").NumberFormat = "General"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E10").NumberFormat = "General"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.42%  '
$ws.Range("E11").NumberFormat = "General"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.010.04'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("E12").NumberFormat = "General"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.47'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.435.35'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("E15").NumberFormat = "General"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.032.15'
$ws.Range("D17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.56'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("E18").NumberFormat = "General"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.96'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.55%  '
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.06'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("E21").NumberFormat = "General"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("E22").NumberFormat = "General"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.44'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("E23").NumberFormat = "General"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.560.22'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("E25").NumberFormat = "General"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("E26").NumberFormat = "General"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.48%  '
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.66'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("E28").NumberFormat = "General"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("E29").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("E30").NumberFormat = "General"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("E31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E32").NumberFormat = "General"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.27'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("E33").NumberFormat = "General"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.68%  '
$ws.Range("E34").NumberFormat = "General"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.34%  '
$ws.Range("E35").NumberFormat = "General"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("E36").NumberFormat = "General"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("E37").NumberFormat = "General"

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("B38").NumberFormat = "General"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C38").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.14'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("E38").NumberFormat = "General"

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Monero'
$ws.Range("B39").NumberFormat = "General"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C39").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '168.52'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E39").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.457.11'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0786'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("E41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.64'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("E42").NumberFormat = "General"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.781'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("E43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.44'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.73%  '
$ws.Range("E44").NumberFormat = "General"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("E45").NumberFormat = "General"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.75%  '
$ws.Range("E46").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.553.41'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("E47").NumberFormat = "General"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("E48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.77'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("E49").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("E50").NumberFormat = "General"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.80%  '
$ws.Range("E51").NumberFormat = "General"
